$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 2943208.11
$ws.Range("C7").Value = -33.75745304248669
$ws.Range("D7").Value = 2985
$ws.Range("E7").Value = 2985
$ws.Range("F7").Value = 985.9993668341708
$ws.Range("G7").Value = 5.100402810982607
